# Apply "All link budget up to and including Mars are done" updates.
# Columns: C = Earth 3U CubeSat, D = Moon 12U Cubesat, E = Mars 6U CubeSat
# (F = Venus Explorer, G = Europa imager, H = BIRD example - untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (S/C transmitter power): Moon 8 -> 28, Mars 5 -> 15
$ws.Range("D3").Value = 28
$ws.Range("E3").Value = 15

# Row 4 (Ground Station transmitter power): Moon 400 -> 100, Mars 400 -> 600
$ws.Range("D4").Value = 100
$ws.Range("E4").Value = 600

# Row 7 (downlink freq): Moon 2.2 -> 8.4
$ws.Range("D7").Value = 8.4

# Row 10 (Antenna D ground station): Mars 10 -> 50
$ws.Range("E10").Value = 50

# Row 11 (Orbit altitude): Mars 500 -> 8000
$ws.Range("E11").Value = 8000

# Row 14 (pointing offset angle S/C): Moon 1 -> 0.5, Mars 1 -> 0.5
$ws.Range("D14").Value = 0.5
$ws.Range("E14").Value = 0.5

# Row 17 (Payload pixel size): Mars 0.2 -> 1
$ws.Range("E17").Value = 1

# Row 19 (Payload duty cycle): Moon 0.5 -> 0.3
$ws.Range("D19").Value = 0.3

# Row 20 (Payload downlink time): Moon 6 -> 12
$ws.Range("D20").Value = 12

# Row 21 (Modulation/coding type): Moon 8FSK -> BPSK_Viterbi, Mars 8FSK -> BPSK_Viterbi
$ws.Range("D21").Value = "BPSK_Viterbi"
$ws.Range("E21").Value = "BPSK_Viterbi"

# Selection moved to E5 on the frozen-pane (bottomRight) view
$ws.Range("E5").Select()
